# Updated cryptos list on Tue May  2 18:56:14 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# re-ranks two coin pairs whose relative order flipped (Aptos/Algorand and
# RenderToken/WEMIXToken), updating their Coin/Link/Price/Volume cells too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.819.55"
$ws.Range("E2").Value = "  +2.18%  "
# Row 3
$ws.Range("D3").Value = "1.876.95"
$ws.Range("E3").Value = "  +2.24%  "
# Row 4
$ws.Range("E4").Value = "  +0.27%  "
# Row 5
$ws.Range("D5").Value = "'327.02"
$ws.Range("E5").Value = "  -0.83%  "
# Row 6
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.20%  "
# Row 7
$ws.Range("D7").Value = "'0.4655"
$ws.Range("E7").Value = "  +1.02%  "
# Row 8
$ws.Range("D8").Value = "'0.3930"
$ws.Range("E8").Value = "  +1.95%  "
# Row 9
$ws.Range("D9").Value = "'0.07919"
$ws.Range("E9").Value = "  +0.95%  "
# Row 10
$ws.Range("D10").Value = "'0.9750"
$ws.Range("E10").Value = "  +2.03%  "
# Row 11
$ws.Range("D11").Value = "'22.36"
$ws.Range("E11").Value = "  +2.07%  "
# Row 12
$ws.Range("D12").Value = "1.862.93"
$ws.Range("E12").Value = "  +2.58%  "
# Row 13
$ws.Range("D13").Value = "'5.757"
$ws.Range("E13").Value = "  +1.10%  "
# Row 14
$ws.Range("D14").Value = "'6.954"
$ws.Range("E14").Value = "  +0.84%  "
# Row 15
$ws.Range("D15").Value = "'0.06998"
$ws.Range("E15").Value = "  +1.95%  "
# Row 16
$ws.Range("D16").Value = "'88.67"
$ws.Range("E16").Value = "  +2.21%  "
# Row 17
$ws.Range("E17").Value = "  +0.25%  "
# Row 18
$ws.Range("D18").Value = "'0.00001013"
$ws.Range("E18").Value = "  +1.92%  "
# Row 19
$ws.Range("D19").Value = "'16.99"
$ws.Range("E19").Value = "  +0.53%  "
# Row 20
$ws.Range("E20").Value = "  +0.06%  "
# Row 21
$ws.Range("D21").Value = "28.830.77"
$ws.Range("E21").Value = "  +2.12%  "
# Row 22
$ws.Range("D22").Value = "'5.345"
$ws.Range("E22").Value = "  +0.24%  "
# Row 23
$ws.Range("E23").Value = "  +1.74%  "
# Row 25
$ws.Range("D25").Value = "2.069.07"
$ws.Range("E25").Value = "  +1.27%  "
# Row 26
$ws.Range("D26").Value = "'153.89"
$ws.Range("E26").Value = "  +0.51%  "
# Row 27
$ws.Range("D27").Value = "'19.41"
# Row 28
$ws.Range("D28").Value = "'5.768"
$ws.Range("E28").Value = "  +1.70%  "
# Row 29
$ws.Range("D29").Value = "'2.010"
$ws.Range("E29").Value = "  +1.86%  "
# Row 30
$ws.Range("D30").Value = "'119.71"
$ws.Range("E30").Value = "  +2.52%  "
# Row 31
$ws.Range("D31").Value = "'0.09379"
# Row 32
$ws.Range("D32").Value = "'0.9403"
$ws.Range("E32").Value = "  +0.45%  "
# Row 33
$ws.Range("D33").Value = "'5.329"
$ws.Range("E33").Value = "  +1.15%  "
# Row 34
$ws.Range("D34").Value = "'1.352"
$ws.Range("E34").Value = "  +2.43%  "
# Row 35
$ws.Range("D35").Value = "'3.355"
$ws.Range("E35").Value = "  -2.57%  "
# Row 36
$ws.Range("D36").Value = "'0.05876"
$ws.Range("E36").Value = "  -1.96%  "
# Row 37
$ws.Range("D37").Value = "'0.02121"
$ws.Range("E37").Value = "  -1.20%  "
# Row 38
$ws.Range("D38").Value = "'1.147"
$ws.Range("E38").Value = "  +0.09%  "
# Row 39
$ws.Range("D39").Value = "'7.931"
$ws.Range("E39").Value = "  +4.43%  "
# Row 40
$ws.Range("D40").Value = "'0.5685"
$ws.Range("E40").Value = "  +1.42%  "
# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'9.976"
$ws.Range("E41").Value = "  +0.11%  "
# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1791"
$ws.Range("E42").Value = "  +1.13%  "
# Row 43
$ws.Range("D43").Value = "'0.07241"
$ws.Range("E43").Value = "  +3.30%  "
# Row 44
$ws.Range("D44").Value = "'11.77"
$ws.Range("E44").Value = "  +1.79%  "
# Row 45
$ws.Range("D45").Value = "'0.5332"
$ws.Range("E45").Value = "  +1.25%  "
# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'2.132"
$ws.Range("E46").Value = "  -4.59%  "
# Row 47
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.134"
$ws.Range("E47").Value = "  -9.25%  "
# Row 48
$ws.Range("D48").Value = "'1.854"
$ws.Range("E48").Value = "  +1.39%  "
# Row 49
$ws.Range("D49").Value = "'113.88"
$ws.Range("E49").Value = "  +1.28%  "
# Row 50
$ws.Range("D50").Value = "'2.363"
$ws.Range("E50").Value = "  +1.48%  "
# Row 51
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.26%  "
